# Auto-generated edit script: update cryptos price/volume table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "1.00", "7.60") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values and mangles the literal formatting (trailing zeros, multi-dot strings, etc).
# NumberFormat "@" forces text entry; ClearFormats() afterwards removes the
# leftover "@" / quote-prefix styling so the cell keeps style index 0, matching
# the original (unstyled) data cells.

$ws.Range("D2").Value = "63.911.99"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "3.053.96"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "3.053.43"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("E10").Value = "  +5.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.07"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -11.41%  "
$ws.Range("E12").Value = "  +6.95%  "
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.50%  "
$ws.Range("D15").Value = "3.551.17"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "63.975.80"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Value = "3.054.98"
$ws.Range("E17").Value = "  +2.12%  "
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.678"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.35"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +14.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.16"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.27"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.21%  "
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.79"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +3.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "443.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0806"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  +11.90%  "
$ws.Range("D41").Value = "2.992.26"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +5.30%  "
$ws.Range("E46").Value = "  +7.87%  "
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("E48").Value = "  +3.89%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("B50").Value = "PEPE"
$ws.Range("C50").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D50").Value = "0.0₃0514"
$ws.Range("E50").Value = "  +4.01%  "
$ws.Range("E51").Value = "  +2.50%  "
